$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column at D, shifting the existing D:K data to E:L
$ws.Columns("D:D").Insert()

# New values for the inserted column D (new quarter of financial data)
$numericValues = @{
    7 = 43407
    8 = 216900
    9 = 146400
    10 = 70500
    13 = 0
    14 = 0
    15 = 6300
    17 = 215000
    18 = 1900
    20 = 300
    21 = 8400
    22 = 0
    23 = 2100
    24 = 600
    25 = 0
    26 = 1500
    27 = 1500
    28 = 0
    30 = 0
    31 = 0
    32 = -300
    33 = 1500
    34 = 0
    35 = 1500
    38 = 43407
    41 = 121200
    42 = 0
    44 = 256900
    45 = 23400
    46 = 401400
    47 = 0
    48 = 103900
    49 = 0
    50 = 0
    51 = 0
    52 = 4700
    53 = 0
    54 = 510000
    57 = 109400
    58 = 25600
    59 = 21300
    60 = 156300
    61 = 0
    62 = 25200
    63 = 0
    64 = 0
    65 = 0
    66 = 181500
    68 = 0
    69 = 0
    70 = 0
    71 = 0
    72 = 753000
    73 = 0
    74 = 0
    75 = 0
    76 = 328500
    77 = 0
    80 = 43407
    81 = 1500
    83 = 6300
    84 = 0
    85 = 0
    86 = 0
    87 = 0
    88 = 0
    89 = -10100
    91 = -6000
    92 = 0
    93 = 0
    94 = -5900
    96 = 0
    97 = 0
    98 = 0
    99 = 0
    100 = 17500
    101 = 0
    102 = 1600
}

$naRows = @(12, 29, 43)

$emptyRows = @(11, 16, 19, 39, 40, 55, 56, 67, 82, 90, 95)

# Copy the number formatting from column E (the old column D, now shifted) onto the new column D
foreach ($r in $numericValues.Keys) {
    $ws.Range("E" + $r).Copy()
    $ws.Range("D" + $r).PasteSpecial(-4122)
}
foreach ($r in $naRows) {
    $ws.Range("E" + $r).Copy()
    $ws.Range("D" + $r).PasteSpecial(-4122)
}
foreach ($r in $emptyRows) {
    $ws.Range("E" + $r).Copy()
    $ws.Range("D" + $r).PasteSpecial(-4122)
}

# Write the numeric values
foreach ($r in $numericValues.Keys) {
    $ws.Range("D" + $r).Value = $numericValues[$r]
}

# Write "NA" markers (reuses the existing shared string)
foreach ($r in $naRows) {
    $ws.Range("D" + $r).Value = "NA"
}
